{"js": "// Apply the \"Play Stolen Treasures for Free\" revision:\n//  - shorten the H1 title and drop the separate \"Meta description\" paragraph\n//  - rewrite the four \"What we like\" bullets\n//  - shorten the last \"What we don't like\" bullet and add a new bold\n//    call-to-action paragraph after it\n//  - replace the DALLE image-prompt paragraph with a short CTA sentence\n\nconst body = context.document.body;\n\nasync function replaceExact(findText, replaceText) {\n  const results = body.search(findText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(replaceText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 1) Shorten the H1 title.\nawait replaceExact(\n  \"Play Stolen Treasures for Free - Exciting Treasure Hunt Adventure\",\n  \"Play Stolen Treasures for Free\"\n);\n\n// 2) Remove the whole \"Meta description: ...\" paragraph entirely.\nconst paras1 = body.paragraphs;\nparas1.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < paras1.items.length; i++) {\n  if (paras1.items[i].text.indexOf(\"Meta description\") === 0) {\n    paras1.items[i].delete();\n    await context.sync();\n    break;\n  }\n}\n\n// 3) Rewrite the \"What we like\" bullet list.\nawait replaceExact(\"Engaging theme and stunning graphics\", \"Engaging theme\");\nawait replaceExact(\"Big payouts potential from Hold & Respin feature\", \"Stunning graphics\");\nawait replaceExact(\"Exciting treasure hunt adventure\", \"Hold & Respin feature\");\nawait replaceExact(\"Beautiful graphics and engaging gameplay\", \"Chance to win big\");\n\n// 4) Shorten the \"Slightly lower RTP...\" bullet under \"What we don't like\".\nawait replaceExact(\"Slightly lower RTP than some other slots\", \"Slightly lower RTP\");\n\n// 5) Insert a new bold \"Play Stolen Treasures for Free\" paragraph right\n//    after that bullet (as a plain/Normal paragraph, not a list item).\nconst paras2 = body.paragraphs;\nparas2.load(\"text\");\nawait context.sync();\nlet target = null;\nfor (let i = 0; i < paras2.items.length; i++) {\n  if (paras2.items[i].text.trim() === \"Slightly lower RTP\") {\n    target = paras2.items[i];\n  }\n}\nif (target) {\n  const newPara = target.insertParagraph(\"Play Stolen Treasures for Free\", Word.InsertLocation.after);\n  newPara.style = \"Normal\";\n  newPara.font.bold = true;\n  await context.sync();\n}\n\n// 6) Replace the DALLE image-prompt paragraph text with the new CTA line.\nawait replaceExact(\n  \"Promopt for DALLE: Create a feature image for Stolen Treasures that captures the adventurous spirit of the game. The image should be in a cartoon style and prominently feature a happy Maya warrior wearing glasses. The warrior should be holding a treasure chest filled with gold and jewels, with a backdrop featuring the ruins of an ancient temple or city. The overall tone should be fun and exciting, with bright colors and playful details that showcase the game's thrilling treasure hunt theme. Make sure the image is eye-catching and appealing to casino players who are in search of new and exciting games to try their luck at.\",\n  \"Read our review of Stolen Treasures and play this thrilling slot game for free.\"\n);\n", "ps1": "# Apply the \"Play Stolen Treasures for Free\" revision:\n#  - shorten the H1 title and drop the separate \"Meta description\" paragraph\n#  - rewrite the four \"What we like\" bullets\n#  - shorten the last \"What we don't like\" bullet and add a new bold\n#    call-to-action paragraph after it\n#  - replace the DALLE image-prompt paragraph with a short CTA sentence\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($findText, $replaceText) {\n    $d.Content.Find.Execute(\n        $findText, $true, $false, $false, $false, $false,\n        $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 1) Shorten the H1 title.\nReplace-ExactText \"Play Stolen Treasures for Free - Exciting Treasure Hunt Adventure\" \"Play Stolen Treasures for Free\"\n\n# 2) Remove the whole \"Meta description: ...\" paragraph entirely.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7).StartsWith(\"Meta description\")) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 3) Rewrite the \"What we like\" bullet list.\nReplace-ExactText \"Engaging theme and stunning graphics\" \"Engaging theme\"\nReplace-ExactText \"Big payouts potential from Hold & Respin feature\" \"Stunning graphics\"\nReplace-ExactText \"Exciting treasure hunt adventure\" \"Hold & Respin feature\"\nReplace-ExactText \"Beautiful graphics and engaging gameplay\" \"Chance to win big\"\n\n# 4) Shorten the \"Slightly lower RTP...\" bullet under \"What we don't like\".\nReplace-ExactText \"Slightly lower RTP than some other slots\" \"Slightly lower RTP\"\n\n# 5) Insert a new bold \"Play Stolen Treasures for Free\" paragraph right\n#    after that bullet (as a plain/Normal paragraph, not a list item).\n$idx = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -ceq \"Slightly lower RTP\") {\n        $idx = $i\n    }\n}\nif ($idx -gt 0) {\n    $target = $d.Paragraphs.Item($idx)\n    $target.Range.InsertParagraphAfter() | Out-Null\n    $newPara = $d.Paragraphs.Item($idx + 1)\n    $newPara.Range.Text = \"Play Stolen Treasures for Free\"\n    $newPara.Style = \"Normal\"\n    $textRange = $newPara.Range.Duplicate()\n    $textRange.MoveEnd(1, -1) | Out-Null\n    $textRange.Font.Bold = 1\n}\n\n# 6) Replace the DALLE image-prompt paragraph text with the new CTA line.\nReplace-ExactText \"Promopt for DALLE: Create a feature image for Stolen Treasures that captures the adventurous spirit of the game. The image should be in a cartoon style and prominently feature a happy Maya warrior wearing glasses. The warrior should be holding a treasure chest filled with gold and jewels, with a backdrop featuring the ruins of an ancient temple or city. The overall tone should be fun and exciting, with bright colors and playful details that showcase the game's thrilling treasure hunt theme. Make sure the image is eye-catching and appealing to casino players who are in search of new and exciting games to try their luck at.\" \"Read our review of Stolen Treasures and play this thrilling slot game for free.\"\n"}
